# Update gh-pages to output generated at 2291077
#
# Adds a new "Cover" column (J) to every sheet, populates the cover image
# URLs for the data rows on the "展览" / "全部类型" sheets, and refreshes
# the "想去人数" (F column) counts that changed since the last scrape.

$wb = $excel.ActiveWorkbook

# Row data shared by the "展览" (sheet 1) and "全部类型" (sheet 4) sheets:
# row number -> new F value (want-to-go count, $null = unchanged) and the
# new J value (cover image URL).
$rowsF = @{
    2  = 406
    3  = 1390
    4  = 6776
    5  = 371
    6  = 196
    7  = 3384
    8  = 29
    9  = 23
    10 = 47
    11 = 841
    13 = 5408
}

$rowsCover = @{
    2  = "//i2.hdslb.com/bfs/openplatform/202312/VBekVPuH1703840712015.jpeg"
    3  = "//i1.hdslb.com/bfs/openplatform/202311/2v00jbxM1698999146733.jpeg"
    4  = "//i1.hdslb.com/bfs/openplatform/202312/iJ1Dnmla1702029064983.jpeg"
    5  = "//i0.hdslb.com/bfs/openplatform/202312/9ClQwbVE1703668101900.jpeg"
    6  = "//i0.hdslb.com/bfs/openplatform/202312/aHzqArm61703662347629.jpeg"
    7  = "//i0.hdslb.com/bfs/openplatform/202312/tBk3WVyX1702968658234.jpeg"
    8  = "//i0.hdslb.com/bfs/openplatform/202401/MSS7qIQp1704695420767.jpeg"
    9  = "//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg"
    10 = "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg"
    11 = "//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg"
    12 = "//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg"
    13 = "//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg"
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New "Cover" header in column J.
    $ws.Cells.Item(1, 10).Value = "Cover"

    foreach ($row in $rowsF.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsF[$row]
    }

    foreach ($row in $rowsCover.Keys) {
        $ws.Cells.Item($row, 10).Value = $rowsCover[$row]
    }
}

# "演出" and "本地生活" sheets only gain the new empty "Cover" header.
foreach ($sheetName in @("演出", "本地生活")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(1, 10).Value = "Cover"
}
